# Weekly update: insert a new price record for Damasco at
# "Macroferia Regional de Talca" as the new first row of the data
# (row 25), pushing the existing rows 25-46 down to rows 26-47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 25 (shifts rows 25..46 down to 26..47).
$ws.Rows.Item(25).Insert()

# Populate the new row 25 with the new record.
$ws.Range("A25").Value = 5
$ws.Range("B25").Value = "Macroferia Regional de Talca"
$ws.Range("C25").Value = "Maule"
$ws.Range("D25").Value = 44907
$ws.Range("E25").Value = 7
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100103
$ws.Range("H25").Value = "Frutos de hueso (carozo)"
$ws.Range("I25").Value = 100103003
$ws.Range("J25").Value = "Damasco"
$ws.Range("K25").Value = "Castle Brite"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 120
$ws.Range("N25").Value = 20000
$ws.Range("O25").Value = 20000
$ws.Range("P25").Value = 20000
$ws.Range("Q25").Value = "`$/caja 18 kilos"
$ws.Range("R25").Value = "Región de O'Higgins"
$ws.Range("S25").Value = 1111
$ws.Range("T25").Value = 18
